$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the date-strings in columns K:P (date fields) for every data row
# from the old dotted format (DD.MM.YYYY) to the ISO format (YYYY-MM-DD),
# each row getting the sequential dates 2021-01-20 .. 2021-01-25.
$dates = @("2021-01-20", "2021-01-21", "2021-01-22", "2021-01-23", "2021-01-24", "2021-01-25")
$cols = @("K", "L", "M", "N", "O", "P")

for ($row = 2; $row -le 10; $row++) {
    for ($i = 0; $i -lt 6; $i++) {
        $ws.Range($cols[$i] + $row).Value = $dates[$i]
    }
}

# Match the author's last-used selection when the file was resaved.
$ws.Range("P7").Select()
